# Apply the Alvearie FHIR IG metadata refresh (StructureDefinition
# wh-payer-substance-abuse-inpatient-coverage-indicator) to the workbook.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# New publication date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to describe a "Contact" (with no displayable ContactDetail);
# it is replaced by a Jurisdiction entry.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row with no useful value - remove it,
# shifting Description/Purpose/... up by one row.
$meta.Rows.Item(11).Delete()

# --- "Elements" sheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The root Extension row's Short/Definition columns now reflect the
# profile's own title/description instead of the generic placeholders.
$elements.Range("K2").Value = "Substance Abuse Inpatient Coverage Indicator"
$elements.Range("L2").Value = "Indicator of Substance Abuse (chemical dependency) inpatient benefit coverage for the member. This finer granularity of MHSA benefit coverage may be used in HEDIS reporting."

# Column K widens to fit the new, longer "Short" text (best-fit recompute).
$elements.Columns.Item(11).ColumnWidth = 42.83
